# Update "想去人数" (interest counts) figures in column F across the
# workbook's sheets, reflecting the regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14062
$ws1.Range("F8").Value = 13929
$ws1.Range("F9").Value = 14951
$ws1.Range("F30").Value = 129
$ws1.Range("F31").Value = 64
$ws1.Range("F32").Value = 322

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 9

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14062
$ws4.Range("F4").Value = 9
$ws4.Range("F9").Value = 13929
$ws4.Range("F10").Value = 14951
$ws4.Range("F32").Value = 129
$ws4.Range("F33").Value = 64
$ws4.Range("F34").Value = 322
